# Upgrade left table: add 2023 column (K) with data, matching existing
# formatting of the table, and extend the custom column-width range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend column widths (B:J were custom 8.71-wide; widen the custom
#        range out to column O, matching the new column count used elsewhere
#        in the workbook family) -----------------------------------------
$refWidth = $ws.Columns.Item(10).ColumnWidth   # width used by column J
for ($i = 11; $i -le 15; $i++) {
    $ws.Columns.Item($i).ColumnWidth = $refWidth
}

# --- 2. Add the "2023" column data, copying the style pattern from the
#        previous last column (J) so number formats / fonts / fills match
#        -------------------------------------------------------------------
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1941
$ws.Range("K5").Value = 401
$ws.Range("K6").Value = 1540

# --- 3. Borders: column J is no longer the right-most column of the table,
#        so drop its old "closing" bottom border on the header row, and give
#        the new right-most column (K) the closing borders instead ---------

# Row 3 (years header): previously every header cell had both a top AND a
# bottom border (the bottom one duplicated row 4's top border as a visual
# "closer"). Now that K3 exists, only the top border remains on B3:J3.
$ws.Range("B3:J3").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none

# K column gets the right-hand closing border all the way down, matching the
# corners: top+right on row3/row4, right-only on row5, right+bottom on row6.
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1       # xlEdgeRight = thin
$ws.Range("K3:K6").Borders.Item(10).Weight = 2          # xlThin

$ws.Range("K3").Borders.Item(9).LineStyle = -4142       # xlEdgeBottom -> none
$ws.Range("K4").Borders.Item(9).LineStyle = -4142       # xlEdgeBottom -> none
